$wb = $excel.ActiveWorkbook

# --- PONTOS sheet: update a couple of point values ---
$pontos = $wb.Worksheets.Item("PONTOS")
$pontos.Range("B9").Value = 1
$pontos.Range("C9").Value = 1.5

# Move the active selection on PONTOS to C9
$pontos.Range("C9").Select()

# --- FREQ sheet: fill in new day column N (rows 3-20) with "P" (presente) ---
$freq = $wb.Worksheets.Item("FREQ")
for ($r = 3; $r -le 20; $r++) {
    $freq.Cells.Item($r, 14).Value = "P"
}

# Move the active selection on FREQ to N19 (mirrors the recorded cursor move),
# and leave FREQ as the active sheet (it's tabSelected in both before/after).
$freq.Range("N19").Select()
